$d = $word.ActiveDocument

# Target paragraph 12: the empty paragraph right after the " است." paragraph,
# right before the final sectPr. It currently only contains pPr/rPr/lang.
$p = $d.Paragraphs(12)
$r = $p.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">۶. اضافه کردن ضریب </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>منظم‌سازی</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> در یک </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>دسته‌بند</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> از بیش برازش آن جلوگیری کرده و تاثیر </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>نویز</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> را روی آن کمتر </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>می‌کند</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>هرچه</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> ضر</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ب</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>منظم‌سازی</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> در توابع هز</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>نه</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> ب</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>شتر</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> شود، اطم</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>نان</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> ب</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ش</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> از حد کاهش </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>افته</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> و مقدار احتمال</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> که توسط مدل </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>براي</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> برچسب خروج</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> تول</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>د</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> م</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>شود</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> ن</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ز</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> کاهش </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>م</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>‌</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ی</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ابد</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:bookmarkStart w:id="1" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="1"/>
    </w:p>'

$r.InsertXML($xml)
